$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update columns B:D for rows 2:25
$arrBD = New-Object 'object[,]' 24,3
$arrBD[0,0] = 11.96456286565524
$arrBD[0,1] = 8.801405357756863
$arrBD[0,2] = 3.772924485062428
$arrBD[1,0] = 11.34279106816993
$arrBD[1,1] = 8.348062187880306
$arrBD[1,2] = 3.69361935843098
$arrBD[2,0] = 10.94323668414463
$arrBD[2,1] = 8.055207346206162
$arrBD[2,2] = 3.643533991449324
$arrBD[3,0] = 10.77610577943572
$arrBD[3,1] = 7.9322965003129
$arrBD[3,2] = 3.62279049056314
$arrBD[4,0] = 10.74809889869187
$arrBD[4,1] = 7.911674145520613
$arrBD[4,2] = 3.619326421313865
$arrBD[5,0] = 10.94099990941261
$arrBD[5,1] = 8.053564071208964
$arrBD[5,2] = 3.643255564238816
$arrBD[6,0] = 11.75395571061043
$arrBD[6,1] = 8.648143320417944
$arrBD[6,2] = 3.745878626054751
$arrBD[7,0] = 13.2012850025827
$arrBD[7,1] = 9.696413105269441
$arrBD[7,2] = 3.93539630774414
$arrBD[8,0] = 14.16898751053592
$arrBD[8,1] = 10.39209109142616
$arrBD[8,2] = 4.0666402016398
$arrBD[9,0] = 14.58753690874995
$arrBD[9,1] = 10.69200778891165
$arrBD[9,2] = 4.124450742821742
$arrBD[10,0] = 14.7428584030605
$arrBD[10,1] = 10.80317484338449
$arrBD[10,2] = 4.146057792850065
$arrBD[11,0] = 14.70954912011616
$arrBD[11,1] = 10.77934032274683
$arrBD[11,2] = 4.141417169754708
$arrBD[12,0] = 14.60037920988792
$arrBD[12,1] = 10.70120188909038
$arrBD[12,2] = 4.12623413789642
$arrBD[13,0] = 14.53309458532749
$arrBD[13,1] = 10.65302601468076
$arrBD[13,2] = 4.116896685981076
$arrBD[14,0] = 14.14119308480746
$arrBD[14,1] = 10.37215552120118
$arrBD[14,2] = 4.062822978798433
$arrBD[15,0] = 13.89517982449455
$arrBD[15,1] = 10.19559188181214
$arrBD[15,2] = 4.029156593125735
$arrBD[16,0] = 13.75164311324017
$arrBD[16,1] = 10.09248038739698
$arrBD[16,2] = 4.009615080313037
$arrBD[17,0] = 13.70269629501639
$arrBD[17,1] = 10.05730208997586
$arrBD[17,2] = 4.002968566034157
$arrBD[18,0] = 13.92157958008825
$arrBD[18,1] = 10.21454870051659
$arrBD[18,2] = 4.032758900163428
$arrBD[19,0] = 14.63253161278786
$arrBD[19,1] = 10.72421847901659
$arrBD[19,2] = 4.130701578570988
$arrBD[20,0] = 15.07865960344122
$arrBD[20,1] = 11.04329194220222
$arrBD[20,2] = 4.193049161486117
$arrBD[21,0] = 14.84226376835286
$arrBD[21,1] = 10.87428651254837
$arrBD[21,2] = 4.159929137368766
$arrBD[22,0] = 13.90965078830017
$arrBD[22,1] = 10.20598331294688
$arrBD[22,2] = 4.031130876363633
$arrBD[23,0] = 12.82617718540578
$arrBD[23,1] = 9.425749559113413
$arrBD[23,2] = 3.885473232748999
$ws.Range("B2:D25").Value = $arrBD

# Update columns F:G for rows 2:25
$arrFG = New-Object 'object[,]' 24,2
$arrFG[0,0] = 18.86780680257212
$arrFG[0,1] = 3.591024579417727
$arrFG[1,0] = 18.8513307742994
$arrFG[1,1] = 3.592772287137118
$arrFG[2,0] = 18.84812358596501
$arrFG[2,1] = 3.593902828334248
$arrFG[3,0] = 18.84855428127793
$arrFG[3,1] = 3.594378021260976
$arrFG[4,0] = 18.84873071918852
$arrFG[4,1] = 3.594457803033553
$arrFG[5,0] = 18.84812235981918
$arrFG[5,1] = 3.593909178231333
$arrFG[6,0] = 18.86069261116054
$arrFG[6,1] = 3.591615293600293
$arrFG[7,0] = 18.9400644039153
$arrFG[7,1] = 3.587570754520698
$arrFG[8,0] = 19.031477752562
$arrFG[8,1] = 3.584873039449904
$arrFG[9,0] = 19.08015457678044
$arrFG[9,1] = 3.583704634727616
$arrFG[10,0] = 19.0995959163563
$arrFG[10,1] = 3.583270600262051
$arrFG[11,0] = 19.09536422516524
$arrFG[11,1] = 3.583363703726219
$arrFG[12,0] = 19.08173388803282
$arrFG[12,1] = 3.583668758030424
$arrFG[13,0] = 19.07351587642425
$arrFG[13,1] = 3.583856707298988
$arrFG[14,0] = 19.02843839451217
$arrFG[14,1] = 3.584950577193644
$arrFG[15,0] = 19.00259385792284
$arrFG[15,1] = 3.58563666214272
$arrFG[16,0] = 18.9883971428198
$arrFG[16,1] = 3.586036816920589
$arrFG[17,0] = 18.98370548979486
$arrFG[17,1] = 3.586173254645937
$arrFG[18,0] = 19.00527595479451
$arrFG[18,1] = 3.585563054477069
$arrFG[19,0] = 19.0857101787791
$arrFG[19,1] = 3.583578928123359
$arrFG[20,0] = 19.14415018959545
$arrFG[20,1] = 3.582331217744282
$arrFG[21,0] = 19.11242665776408
$arrFG[21,1] = 3.582992671132424
$arrFG[22,0] = 19.00406131696288
$arrFG[22,1] = 3.585596314696526
$arrFG[23,0] = 18.91275307256808
$arrFG[23,1] = 3.588616623684047
$ws.Range("F2:G25").Value = $arrFG

# Update columns I:I for rows 2:25
$arrII = New-Object 'object[,]' 24,1
$arrII[0,0] = 16.08450991216507
$arrII[1,0] = 16.19128804063394
$arrII[2,0] = 16.26095388935468
$arrII[3,0] = 16.29037431790672
$arrII[4,0] = 16.29532180547576
$arrII[5,0] = 16.26134649003671
$arrII[6,0] = 16.12047476475101
$arrII[7,0] = 15.87682606079939
$arrII[8,0] = 15.71774765006434
$arrII[9,0] = 15.64972037263279
$arrII[10,0] = 15.62458552413423
$arrII[11,0] = 15.62997091031962
$arrII[12,0] = 15.64763996498512
$arrII[13,0] = 15.65854429754055
$arrII[14,0] = 15.72228083855697
$arrII[15,0] = 15.76249350460694
$arrII[16,0] = 15.7860309949015
$arrII[17,0] = 15.79407044625321
$arrII[18,0] = 15.75817053509393
$arrII[19,0] = 15.64243314200532
$arrII[20,0] = 15.57043946335202
$arrII[21,0] = 15.60852949733239
$arrII[22,0] = 15.76012364553277
$arrII[23,0] = 15.93924395004476
$ws.Range("I2:I25").Value = $arrII

# Update columns N:O for rows 2:25
$arrNO = New-Object 'object[,]' 24,2
$arrNO[0,0] = 15.66766933099571
$arrNO[0,1] = 16.80014617381049
$arrNO[1,0] = 15.68426547690728
$arrNO[1,1] = 16.8508199938014
$arrNO[2,0] = 15.69617869804967
$arrNO[2,1] = 16.88709242250533
$arrNO[3,0] = 15.70146762440504
$arrNO[3,1] = 16.90316548451424
$arrNO[4,0] = 15.70237209642761
$arrNO[4,1] = 16.90591225798765
$arrNO[5,0] = 15.69624826716918
$arrNO[5,1] = 16.88730396659476
$arrNO[6,0] = 15.67303459553761
$arrNO[6,1] = 16.81654472918864
$arrNO[7,0] = 15.64114517519807
$arrNO[7,1] = 16.71894093948198
$arrNO[8,0] = 15.62597316410574
$arrNO[8,1] = 16.67260192347664
$arrNO[9,0] = 15.62085194050762
$arrNO[9,1] = 16.6570810157421
$arrNO[10,0] = 15.6191676204987
$arrNO[10,1] = 16.65200640356145
$arrNO[11,0] = 15.61951904632764
$arrNO[11,1] = 16.65306355719891
$arrNO[12,0] = 15.62070826667236
$arrNO[12,1] = 16.65664741677483
$arrNO[13,0] = 15.62146987137397
$arrNO[13,1] = 16.6589472783427
$arrNO[14,0] = 15.62634358766493
$arrNO[14,1] = 16.67372841826354
$arrNO[15,0] = 15.62978875976683
$arrNO[15,1] = 16.68422256721406
$arrNO[16,0] = 15.6319379968885
$arrNO[16,1] = 16.69078158433843
$arrNO[17,0] = 15.63269452030081
$arrNO[17,1] = 16.69309207729643
$arrNO[18,0] = 15.62940467022035
$arrNO[18,1] = 16.68305128571549
$arrNO[19,0] = 15.62035205259906
$arrNO[19,1] = 16.65557293719468
$arrNO[20,0] = 15.61592131317844
$arrNO[20,1] = 16.64229491415621
$arrNO[21,0] = 15.61815050834315
$arrNO[21,1] = 16.64895236110834
$arrNO[22,0] = 15.6295777919634
$arrNO[22,1] = 16.68357918454212
$arrNO[23,0] = 15.64831804019256
$arrNO[23,1] = 16.7409086375336
$ws.Range("N2:O25").Value = $arrNO
